$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: set a cell to a NUMBER while reusing the number-style of a sibling
# cell that already carries the right style (xlPasteFormats after the value
# write so the value-coercion doesn't get clobbered by the paste).
# ---------------------------------------------------------------------------
function Set-NumberLike($targetAddr, $value, $styleSourceAddr) {
    $ws.Range($targetAddr).Value = $value
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Helper: set a cell to TEXT (shared-string) while reusing the text-style of
# a sibling cell. NumberFormat "@" forces the Value write to be stored as
# text instead of being silently re-parsed as a number; the subsequent
# PasteFormats restores the real target style (e.g. s="14").
# ---------------------------------------------------------------------------
function Set-TextLike($targetAddr, $value, $styleSourceAddr) {
    $ws.Range($targetAddr).NumberFormat = "@"
    $ws.Range($targetAddr).Value = $value
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($targetAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Title block: "Volume 30   Number  2" -> "...Number  3"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  3"

# Report week header: 1/9/2023 - 1/15/2023 -> 1/16/2023 - 1/22/2023
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 3
$ws.Range("J15").Value = 2
Set-NumberLike "L15" -100 "K15"
Set-NumberLike "M15" -100 "K15"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 57.142857142857
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -11.764705882352
$ws.Range("N16").Value = -85.576923076923

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 28
$ws.Range("H17").Value = 115.384615384615
$ws.Range("I17").Value = 22
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 100
$ws.Range("L17").Value = 69.230769230769
$ws.Range("M17").Value = 83.333333333333
$ws.Range("N17").Value = -53.191489361702

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -62.5
$ws.Range("N18").Value = -95.348837209302

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = -35.849056603773
$ws.Range("I19").Value = 26
$ws.Range("J19").Value = 41
$ws.Range("K19").Value = -36.585365853658
$ws.Range("L19").Value = -33.333333333333
$ws.Range("M19").Value = 36.842105263157
$ws.Range("N19").Value = -61.194029850746

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -62.5
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 13
$ws.Range("J20").Value = 32
$ws.Range("K20").Value = -59.375
$ws.Range("L20").Value = 62.5
$ws.Range("M20").Value = 85.714285714285
$ws.Range("N20").Value = -90.441176470588

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -13.333333333333
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 128
$ws.Range("H21").Value = -12.5
$ws.Range("I21").Value = 82
$ws.Range("J21").Value = 104
$ws.Range("K21").Value = -21.153846153846
$ws.Range("L21").Value = 5.128205128205
$ws.Range("M21").Value = 13.888888888888
$ws.Range("N21").Value = -83.501006036217

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-TextLike "D22" "0" "C22"
Set-TextLike "E22" "***.*" "N22"
$ws.Range("M22").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 4
$ws.Range("K23").Value = 100
$ws.Range("L23").Value = 33.333333333333
$ws.Range("M23").Value = 300

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -32
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 40.963855421686
$ws.Range("I24").Value = 90
$ws.Range("J24").Value = 57
$ws.Range("K24").Value = 57.894736842105
$ws.Range("L24").Value = 95.652173913043
$ws.Range("M24").Value = 210.344827586207

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -17.391304347826
$ws.Range("I25").Value = 17
$ws.Range("J25").Value = 20
$ws.Range("K25").Value = -15
$ws.Range("L25").Value = -22.727272727272
$ws.Range("M25").Value = -43.333333333333

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("J26").Value = 3
Set-NumberLike "L26" -100 "K26"

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = -40

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic.
# ---------------------------------------------------------------------------
Set-TextLike "D28" "0" "C28"
Set-TextLike "E28" "***.*" "L28"
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -33.333333333333

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc.
# ---------------------------------------------------------------------------
Set-TextLike "D29" "0" "C29"
Set-TextLike "E29" "***.*" "L29"
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = -33.333333333333
